$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 and A3 with the combined card-data tuples
$ws.Range("A2").Value = '(''Hangarback Walker'', [''{X}{X}'', ''Artifact Creature — Construct'', ''Hangarback Walker enters the battlefield with X +1/+1 counters on it.'', ''When Hangarback Walker dies, create a 1/1 colorless Thopter artifact creature token with flying for each +1/+1 counter on Hangarback Walker.'', ''{1}, {T}: Put a +1/+1 counter on Hangarback Walker.'', ''0/0''])'
$ws.Range("A3").Value = '(''Reliquary Tower'', [''Land'', ''You have no maximum hand size.'', ''{T}: Add {C}.''])'

# Remove the now-obsolete rows 4 through 12 (their content was folded into A2/A3)
$ws.Range("A4:A12").EntireRow.Delete()
